# edit.ps1 - applies the "create readme, edit manual (#79)" changes:
#  1) Bump the cached date shown by the two title-slide date fields
#     (slide layouts "1_Titelfolie" and "3_Titelfolie") from 30.07.2025
#     to 31.07.2025.
#  2) On slide 2, shorten the "Button 1" description inside the grouped
#     callout textbox from "Pauses and resumes the game" to
#     "Pause / Resume" (the textbox auto-shrinks to fit, which also
#     updates its height).

$p = $ppt.ActivePresentation

# --- 1) Update the two date fields living on the slide master's custom layouts ---
$master = $p.SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "30.07.2025") {
                $tr.Text = "31.07.2025"
            }
        }
    }
}

# --- 2) Update the "Button 1" helper text on slide 2 ---
$slide2 = $p.Slides.Item(2)

function Find-ShapeByName($shapes, $name) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Name -eq $name) {
            return $candidate
        }
    }
    return $null
}

$group = Find-ShapeByName $slide2.Shapes "Gruppieren 6"
$textBox = Find-ShapeByName $group.GroupItems "Textfeld 4"

$tr = $textBox.TextFrame.TextRange
$oldLine = "Pauses and resumes the game"
$newLine = "Pause / Resume"
$start = $tr.Text.IndexOf($oldLine)
if ($start -ge 0) {
    $sub = $tr.Characters($start + 1, $oldLine.Length)
    $sub.Text = $newLine
}
